$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from D1 into E1:G1 (bold, centered, bordered), then set header text
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "Complexidade do projeto (0-100)"
$ws.Range("F1").Value = "Demanda do mercado (0-100)"
$ws.Range("G1").Value = "Dependência da tecnologia de terceiros (0-100)"

# New data for columns E (Complexidade), F (Demanda), G (Dependencia) - rows 2..101
$data = @(
    @(31,23,60),
    @(8,31,63),
    @(98,46,62),
    @(18,85,68),
    @(47,22,21),
    @(79,65,92),
    @(2,26,66),
    @(19,1,75),
    @(23,89,25),
    @(53,16,15),
    @(32,32,50),
    @(23,8,100),
    @(74,42,85),
    @(71,47,56),
    @(35,38,28),
    @(37,92,77),
    @(83,41,91),
    @(98,25,68),
    @(88,98,46),
    @(98,49,93),
    @(24,24,61),
    @(92,23,68),
    @(17,12,75),
    @(81,59,15),
    @(65,6,89),
    @(53,56,89),
    @(34,35,47),
    @(79,44,84),
    @(60,19,38),
    @(40,64,99),
    @(99,7,32),
    @(32,15,93),
    @(67,13,100),
    @(32,75,22),
    @(13,86,9),
    @(20,14,68),
    @(47,91,99),
    @(19,97,33),
    @(7,65,51),
    @(6,31,94),
    @(66,86,9),
    @(16,62,18),
    @(32,85,57),
    @(47,50,95),
    @(75,24,0),
    @(58,57,68),
    @(85,62,3),
    @(21,61,15),
    @(29,21,23),
    @(37,57,79),
    @(50,57,1),
    @(53,85,91),
    @(7,48,31),
    @(26,51,90),
    @(26,41,83),
    @(97,69,23),
    @(20,14,11),
    @(29,53,49),
    @(96,59,34),
    @(27,100,32),
    @(63,96,32),
    @(96,7,60),
    @(68,52,50),
    @(60,59,42),
    @(47,4,100),
    @(18,67,11),
    @(3,5,66),
    @(34,95,64),
    @(63,93,32),
    @(48,46,39),
    @(16,98,73),
    @(43,54,42),
    @(91,39,43),
    @(29,51,28),
    @(92,15,12),
    @(45,12,11),
    @(5,29,94),
    @(98,18,45),
    @(36,16,1),
    @(23,62,34),
    @(92,18,86),
    @(45,91,80),
    @(52,57,89),
    @(94,54,7),
    @(98,89,92),
    @(59,100,25),
    @(96,89,73),
    @(62,61,89),
    @(84,22,33),
    @(31,8,6),
    @(86,11,67),
    @(32,0,57),
    @(66,57,74),
    @(17,0,28),
    @(24,33,35),
    @(94,95,88),
    @(53,47,20),
    @(57,88,35),
    @(66,0,9),
    @(45,15,100)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = 2 + $i
    $ws.Cells.Item($r, 5).Value = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
}
